# update hotel reviews data
# hotel_info!G2 = English_Reviews_num, hotel_info!H2 = Local_Rank
# These were previously blank text cells; fill them in as text values
# (matching the text typing of every other populated column in this row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# Leading apostrophe forces these numeric-looking values to be stored as
# text (shared strings), consistent with the rest of the row/sheet.
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = "'30"

# Drop the quote-prefix cell style Excel applies for text-forced numbers so
# the cells keep the sheet's default formatting.
$ws.Range("G2:H2").Style = "Normal"
